$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.34640275735485204,
    -0.0099999993219412886,
    -0.0089999993110687626,
    -0.011999999809692241,
    -0.0059999993139818741,
    -0.0059999992857093787,
    -0.019999999161852244,
    -0.019999999152794601,
    -0.0059999992631274424,
    -0.0059999992537100866,
    -0.0044999992663541377,
    -0.0059999992512662637,
    -0.005999999242183307,
    -0.011999999187786159,
    -0.0059999992379156097,
    -0.0059999992366932542,
    -0.0059999992357449017,
    0.0010674999390012374,
    -0.0089999993299367809,
    -0.053028105020640481,
    -0.0089999993128024869,
    -0.0089999993121665511,
    -0.0089999993095286612,
    0.026461655237249637,
    -0.041999999017933831,
    -0.0059999992832757698,
    -0.0068804545131380124,
    -0.0059999992654242718,
    -0.011999999203823108,
    -0.017248229732540921,
    -0.014999999165999256,
    -0.020999999112985002,
    -0.0059999992421690962
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
